$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update a handful of tracked hours / estimates in the backlog ---
$ws.Range("C5").Value = 3
$ws.Range("D5").Value = 9
$ws.Range("C6").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("C12").Value = 3
$ws.Range("F12").Value = 3

# --- Capture the chart's anchor height (in points) before we delete a row,
#     so we can shrink the chart by exactly one row afterwards (mirrors
#     Excel's "move and size with cells" behaviour for a twoCellAnchor) ---
$co = $ws.ChartObjects().Item(1)
$rowHeightPts = $ws.Rows(18).Height

# --- Remove the "remove a meal from planned meals" backlog row entirely;
#     everything below shifts up by one row ---
$ws.Rows("18:18").Delete()

# --- Keep the chart object's footprint anchored correctly now that the
#     sheet is one row shorter ---
$co2 = $ws.ChartObjects().Item(1)
$co2.Height = $co2.Height - $rowHeightPts

# --- Fix up the (now shifted) filter database range ---
$fdb = $wb.Names.Item("Sheet1!_FilterDatabase")
$fdb.RefersTo = "=Sheet1!`$B`$25:`$G`$28"

# --- Re-point the burndown chart series at the new Estimate Totals row ---
$chart = $co2.Chart
$ser = $chart.SeriesCollection().Item(1)
$ser.Formula = "=SERIES(,,Sheet1!`$F`$25:`$G`$25,1)"

# --- Restore the last active selection ---
$ws.Range("H16").Select()
